$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Table 2 Authorizations": rename the OIMP field header, update the
#    syntax-check comment left on A1, and move the active selection to P7.
# ---------------------------------------------------------------------------
$wsAuth = $wb.Worksheets.Item("Table 2 Authorizations")

$wsAuth.Range("P7").Value = "OIMP authorized by the Party2"

$wsAuth.Range("A1").Comment.Text("The field 'OIMP authorized by the Party' cannot be found in worksheet")

$wsAuth.Activate()
$wsAuth.Range("P7").Select()

# ---------------------------------------------------------------------------
# 2. "Syntax check results": a new error line is reported for Table 2
#    Authorizations, so insert a row for it right after the existing
#    Table 2 Authorizations row (row 16), push everything else down, and
#    link it back to the offending sheet.
# ---------------------------------------------------------------------------
$wsSyntax = $wb.Worksheets.Item("Syntax check results")

$wsSyntax.Rows.Item(17).Insert()

$wsSyntax.Range("D17").Value = "Link"
$wsSyntax.Range("E17").Value = "The field 'OIMP authorized by the Party' cannot be found in worksheet"
$wsSyntax.Hyperlinks.Add($wsSyntax.Range("D17"), "#'Table 2 Authorizations'!A1")

# Formatting touch-up: bold the file name, italicise the "Structure check"
# caption and give every "Link" cell the familiar blue/underlined look.
$wsSyntax.Range("A1").Font.Bold = $true
$wsSyntax.Range("B2").Font.Italic = $true

$linkCells = @("D6", "D8", "D10", "D12", "D14", "D16", "D17", "D19", "D21", "D23")
foreach ($cellRef in $linkCells) {
    $cell = $wsSyntax.Range($cellRef)
    $cell.Font.Underline = $true
    $cell.Font.Color = 16711680
}
